$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (testing person 1 / invitation 60107): RSVP flipped Maybe -> No,
# date opened refreshed, and a diet info ("Vegan") recorded.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "No"
$ws.Range("H2").Value = "2015-04-23 10:11"
$ws.Range("I2").Value = "Vegan"

# ---------------------------------------------------------------------------
# Row 3 (testing person 2 / invitation 60107): RSVP flipped Maybe -> Yes,
# now coming (# coming 0 -> 1), date opened refreshed, diet info GlutenFree.
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "Yes"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "2015-04-23 10:11"
$ws.Range("I3").Value = "GlutenFree"

# ---------------------------------------------------------------------------
# Row 4 (the Moskovitzes / invitation 42652): date opened refreshed only.
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = "2015-04-23 10:08"

# ---------------------------------------------------------------------------
# Row 5 (army friend / invitation 20349): date opened refreshed, diet info
# changed from Vegan to GlutenFree.
# ---------------------------------------------------------------------------
$ws.Range("H5").Value = "2015-04-23 10:19"
$ws.Range("I5").Value = "GlutenFree"

# ---------------------------------------------------------------------------
# Row 6 (Guest / invitation 20349): date opened refreshed only.
# ---------------------------------------------------------------------------
$ws.Range("H6").Value = "2015-04-23 10:19"

# ---------------------------------------------------------------------------
# New row 7: a brand-new invitation (88850), mirror the borders/formatting
# used by the other data rows before writing values.
# ---------------------------------------------------------------------------
$ws.Range("A6:M6").Copy()
$ws.Range("A7:M7").PasteSpecial(-4122)

$ws.Range("A7").Value = "88850"
$ws.Range("B7").Value = "sththrtger"
$ws.Range("C7").Value = "srhsth"
$ws.Range("D7").Value = "httshtrh"
$ws.Range("E7").Value = "Maybe"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = "2015-04-21 16:57"
$ws.Range("J7").Value = "Both"
$ws.Range("K7").Value = "Family"
$ws.Range("M7").Value = "avichaidevora.com/invitation/88850"

# ---------------------------------------------------------------------------
# Column width tweaks: D widens to match the "wide" columns, G narrows back
# to the "narrow" width, and H widens to match the "wide" columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
